# Insert two new data rows at row 184 (pushing existing rows 184-248 down to 186-250)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A184:A185").EntireRow.Insert()

# New row 184: Choclero / Primera, fecha 2022-12-29, Región de O'Higgins
$ws.Cells.Item(184, 1).Value = 7
$ws.Cells.Item(184, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(184, 3).Value = "Ñuble"
$ws.Cells.Item(184, 4).Value = 44924
$ws.Cells.Item(184, 5).Value = 16
$ws.Cells.Item(184, 6).Value = 100112024
$ws.Cells.Item(184, 7).Value = "Choclo"
$ws.Cells.Item(184, 8).Value = "Choclero"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 20000
$ws.Cells.Item(184, 11).Value = 250
$ws.Cells.Item(184, 12).Value = 300
$ws.Cells.Item(184, 13).Value = 275
$ws.Cells.Item(184, 14).Value = "`$/unidad"
$ws.Cells.Item(184, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(184, 16).Value = 275
$ws.Cells.Item(184, 17).Value = 1
$ws.Cells.Item(184, 18).Value = "Hortaliza"

# New row 185: Choclero / Segunda, fecha 2022-12-29, Región de O'Higgins
$ws.Cells.Item(185, 1).Value = 7
$ws.Cells.Item(185, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(185, 3).Value = "Ñuble"
$ws.Cells.Item(185, 4).Value = 44924
$ws.Cells.Item(185, 5).Value = 16
$ws.Cells.Item(185, 6).Value = 100112024
$ws.Cells.Item(185, 7).Value = "Choclo"
$ws.Cells.Item(185, 8).Value = "Choclero"
$ws.Cells.Item(185, 9).Value = "Segunda"
$ws.Cells.Item(185, 10).Value = 15000
$ws.Cells.Item(185, 11).Value = 200
$ws.Cells.Item(185, 12).Value = 200
$ws.Cells.Item(185, 13).Value = 200
$ws.Cells.Item(185, 14).Value = "`$/unidad"
$ws.Cells.Item(185, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(185, 16).Value = 200
$ws.Cells.Item(185, 17).Value = 1
$ws.Cells.Item(185, 18).Value = "Hortaliza"
